$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Status: draft -> active
$ws.Cells.Item(6, 2).Value = "active"

# Date: publish date updated
$ws.Cells.Item(8, 2).Value = "2024-12-16T14:50:05-03:00"

# Case Sensitive: false -> true (write via a formula + paste-values so the
# literal word "true" lands as text, not an auto-coerced Boolean)
$caseSensitiveCell = $ws.Cells.Item(17, 2)
$caseSensitiveCell.Formula = "=""true"""
$caseSensitiveCell.Copy()
$caseSensitiveCell.PasteSpecial(-4163)
